$wb = $excel.ActiveWorkbook

# --- Sheet "Joe": bump every year by one and drop the now-superfluous last row ---
$joe = $wb.Worksheets.Item("Joe")
for ($r = 2; $r -le 36; $r++) {
    $cell = $joe.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}
$joe.Rows.Item(37).Delete() | Out-Null
$joe.Range("B3").Select() | Out-Null

# --- Sheet "Debts": update selection ---
$debts = $wb.Worksheets.Item("Debts")
$debts.Range("D1:D1048576").Select() | Out-Null

# --- Sheet "Fixed Assets": insert a new "year" column before the existing column D ---
$assets = $wb.Worksheets.Item("Fixed Assets")
$assets.Columns.Item(4).Insert() | Out-Null
$assets.Range("D1").Value = "year"

# Make "Fixed Assets" the active tab/sheet with column D selected
$assets.Activate() | Out-Null
$assets.Range("D1:D1048576").Select() | Out-Null
